# Applies the "stuff at the bottom of the sheets" commit:
#  - adds a "pair_kind" value ("generic") to column J for the four
#    practice-pair rows (2-5)
#  - appends a new "stim details" block starting at row 27

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- fill in the new J column (pair_kind) for the practice rows ---
$ws.Range("J2").Value = "generic"
$ws.Range("J3").Value = "generic"
$ws.Range("J4").Value = "generic"
$ws.Range("J5").Value = "generic"

# --- new "stim details" block appended below the existing data ---
$ws.Range("A27").Value = "stim details"

$ws.Range("A28").Value = "month"
$ws.Range("B28").Value = "word_type"
$ws.Range("C28").Value = "need_audio"
$ws.Range("D28").Value = "need_image"
$ws.Range("E28").Value = "word"
$ws.Range("F28").Value = "count"
$ws.Range("G28").Value = "find images"

$ws.Cells.Item(29, 1).Value = 6
$ws.Cells.Item(29, 2).Value = "video"

$ws.Cells.Item(30, 1).Value = 6
$ws.Cells.Item(30, 2).Value = "video"

$ws.Cells.Item(31, 1).Value = 7
$ws.Cells.Item(31, 2).Value = "video"

$ws.Cells.Item(32, 1).Value = 7
$ws.Cells.Item(32, 2).Value = "video"

$ws.Cells.Item(33, 1).Value = 6
$ws.Cells.Item(33, 2).Value = "audio"

$ws.Cells.Item(34, 1).Value = 6
$ws.Cells.Item(34, 2).Value = "audio"

$ws.Cells.Item(35, 1).Value = 7
$ws.Cells.Item(35, 2).Value = "audio"

$ws.Cells.Item(36, 1).Value = 7
$ws.Cells.Item(36, 2).Value = "audio"
